# Update scripts with new TPM values for Gdf6-Bmpr1b.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster becomes "FAPs" (was "ECs"), with refreshed TPM-derived metrics.
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7708836666666666
$ws.Range("N2").Value = 2.312651
$ws.Range("O2").Value = 0.5975056510655317
$ws.Range("P2").Value = 0.6900916471389698
$ws.Range("Q2").Value = 0.1790400442343333
$ws.Range("R2").Value = 1.611360398109
$ws.Range("S2").Value = 0.5975056510655317
$ws.Range("T2").Value = 0.6900916471389698

# Row 3: Target cluster becomes "MuSCs" (was "FAPs"), with refreshed TPM-derived metrics.
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.519286
$ws.Range("N3").Value = 1.038572
$ws.Range("O3").Value = 0.4024943489344683
$ws.Range("P3").Value = 0.3099083528610301
$ws.Range("Q3").Value = 0.120605731358
$ws.Range("R3").Value = 0.723634388148
$ws.Range("S3").Value = 0.4024943489344683
$ws.Range("T3").Value = 0.3099083528610301

# Rows 4 and 5 (ECs/Neutrophils-less entries) are no longer part of the dataset.
$ws.Rows("4:5").Delete()
